$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Outcomes")
$ws2 = $wb.Worksheets.Item("Budget & coverage")

# --- Sheet "Outcomes": update recalculated scenario values (rows 18-31) ---
# Row 18
$ws1.Cells.Item(18, 4).Value = 412047.5320061742
$ws1.Cells.Item(18, 5).Value = 413735.1222047611
$ws1.Cells.Item(18, 6).Value = 415679.6962505301
$ws1.Cells.Item(18, 7).Value = 418132.3066504938
$ws1.Cells.Item(18, 8).Value = 421650.610904283
$ws1.Cells.Item(18, 9).Value = 426331.7709152632
$ws1.Cells.Item(18, 10).Value = 432162.9547274425
$ws1.Cells.Item(18, 11).Value = 438896.3488175675
$ws1.Cells.Item(18, 12).Value = 446327.6181845064
$ws1.Cells.Item(18, 13).Value = 454470.9799511385
$ws1.Cells.Item(18, 14).Value = 463135.6302939343
$ws1.Cells.Item(18, 15).Value = 472133.2973482833
$ws1.Cells.Item(18, 16).Value = 481443.2397931197
$ws1.Cells.Item(18, 17).Value = 5696147.108047497

# Row 19
$ws1.Cells.Item(19, 4).Value = 196848.3045727076
$ws1.Cells.Item(19, 5).Value = 199317.4642202598
$ws1.Cells.Item(19, 6).Value = 202157.8764909415
$ws1.Cells.Item(19, 7).Value = 205204.7989916589
$ws1.Cells.Item(19, 8).Value = 208940.8827616725
$ws1.Cells.Item(19, 9).Value = 212955.9865718944
$ws1.Cells.Item(19, 10).Value = 216605.3056433503
$ws1.Cells.Item(19, 11).Value = 221261.5091597556
$ws1.Cells.Item(19, 12).Value = 225718.5271261221
$ws1.Cells.Item(19, 13).Value = 230165.883895225
$ws1.Cells.Item(19, 14).Value = 234648.3851311253
$ws1.Cells.Item(19, 15).Value = 239681.9775283465
$ws1.Cells.Item(19, 16).Value = 244365.0364040733
$ws1.Cells.Item(19, 17).Value = 2837871.938497133

# Row 20
$ws1.Cells.Item(20, 4).Value = 1524003.945075005
$ws1.Cells.Item(20, 5).Value = 1530245.686430316
$ws1.Cells.Item(20, 6).Value = 1537437.911325243
$ws1.Cells.Item(20, 7).Value = 1546509.165574239
$ws1.Cells.Item(20, 8).Value = 1559522.007895598
$ws1.Cells.Item(20, 9).Value = 1576835.802470797
$ws1.Cells.Item(20, 10).Value = 1598403.088873504
$ws1.Cells.Item(20, 11).Value = 1623307.301033701
$ws1.Cells.Item(20, 12).Value = 1650792.683064868
$ws1.Cells.Item(20, 13).Value = 1680911.81858826
$ws1.Cells.Item(20, 14).Value = 1712958.998293167
$ws1.Cells.Item(20, 15).Value = 1746237.877602507
$ws1.Cells.Item(20, 16).Value = 1780671.742417363
$ws1.Cells.Item(20, 17).Value = 21067838.02864457

# Row 21
$ws1.Cells.Item(21, 4).Value = 43270.15530314871
$ws1.Cells.Item(21, 5).Value = 43447.37342562582
$ws1.Cells.Item(21, 6).Value = 43651.57807298033
$ws1.Cells.Item(21, 7).Value = 43909.13290504747
$ws1.Cells.Item(21, 8).Value = 44278.59894875216
$ws1.Cells.Item(21, 9).Value = 44770.17942173857
$ws1.Cells.Item(21, 10).Value = 45382.52680779785
$ws1.Cells.Item(21, 11).Value = 46089.61758099111
$ws1.Cells.Item(21, 12).Value = 46869.99400514474
$ws1.Cells.Item(21, 13).Value = 47725.14905635324
$ws1.Cells.Item(21, 14).Value = 48635.0459416857
$ws1.Cells.Item(21, 15).Value = 49579.9137556322
$ws1.Cells.Item(21, 16).Value = 50557.57439951731
$ws1.Cells.Item(21, 17).Value = 598166.8396244153

# Row 22
$ws1.Cells.Item(22, 4).Value = 393289.4970542707
$ws1.Cells.Item(22, 5).Value = 394900.2614661295
$ws1.Cells.Item(22, 6).Value = 396756.3108029425
$ws1.Cells.Item(22, 7).Value = 399097.2686677324
$ws1.Cells.Item(22, 8).Value = 402455.405783008
$ws1.Cells.Item(22, 9).Value = 406923.4608576441
$ws1.Cells.Item(22, 10).Value = 412489.1860970964
$ws1.Cells.Item(22, 11).Value = 418916.0494307867
$ws1.Cells.Item(22, 12).Value = 426009.0179957981
$ws1.Cells.Item(22, 13).Value = 433781.6616952934
$ws1.Cells.Item(22, 14).Value = 442051.8628511837
$ws1.Cells.Item(22, 15).Value = 450639.9204795036
$ws1.Cells.Item(22, 16).Value = 459526.0374862374
$ws1.Cells.Item(22, 17).Value = 5436835.940667626

# Row 23
$ws1.Cells.Item(23, 4).Value = 0.6852802394633019
$ws1.Cells.Item(23, 5).Value = 0.6846434162517728
$ws1.Cells.Item(23, 6).Value = 0.6837567801787546
$ws1.Cells.Item(23, 7).Value = 0.6830699353360487
$ws1.Cells.Item(23, 8).Value = 0.6823389940398397
$ws1.Cells.Item(23, 9).Value = 0.6817817950919728
$ws1.Cells.Item(23, 10).Value = 0.6816476519808313
$ws1.Cells.Item(23, 11).Value = 0.6811695620910542
$ws1.Cells.Item(23, 12).Value = 0.680989449194164
$ws1.Cells.Item(23, 13).Value = 0.6809350682965977
$ws1.Cells.Item(23, 14).Value = 0.6809363857593569
$ws1.Cells.Item(23, 15).Value = 0.6807510953109484
$ws1.Cells.Item(23, 16).Value = 0.6808120620212049

# Row 24
$ws1.Cells.Item(24, 4).Value = 0.0346234399586491
$ws1.Cells.Item(24, 5).Value = 0.03468647504218673
$ws1.Cells.Item(24, 6).Value = 0.03477215913956811
$ws1.Cells.Item(24, 7).Value = 0.03484658458410779
$ws1.Cells.Item(24, 8).Value = 0.03492602518404566
$ws1.Cells.Item(24, 9).Value = 0.03498891364200997
$ws1.Cells.Item(24, 10).Value = 0.03500994904016883
$ws1.Cells.Item(24, 11).Value = 0.03506006695946892
$ws1.Cells.Item(24, 12).Value = 0.03508238143504402
$ws1.Cells.Item(24, 13).Value = 0.03509121407912972
$ws1.Cells.Item(24, 14).Value = 0.03509278650237659
$ws1.Cells.Item(24, 15).Value = 0.03511110378402864
$ws1.Cells.Item(24, 16).Value = 0.0351068144602926

# Row 25
$ws1.Cells.Item(25, 4).Value = 0.2203012801318656
$ws1.Cells.Item(25, 5).Value = 0.2199661335231438
$ws1.Cells.Item(25, 6).Value = 0.2199285168885287
$ws1.Cells.Item(25, 7).Value = 0.2199795812079206
$ws1.Cells.Item(25, 8).Value = 0.2200102001497022
$ws1.Cells.Item(25, 9).Value = 0.220068871845497
$ws1.Cells.Item(25, 10).Value = 0.2201569611259182
$ws1.Cells.Item(25, 11).Value = 0.2201340913938908
$ws1.Cells.Item(25, 12).Value = 0.2201897672038575
$ws1.Cells.Item(25, 13).Value = 0.220222334856519
$ws1.Cells.Item(25, 14).Value = 0.2202384877928925
$ws1.Cells.Item(25, 15).Value = 0.2202176520111073
$ws1.Cells.Item(25, 16).Value = 0.2202538041817151

# Row 31
$ws1.Cells.Item(31, 4).Value = 91.55735096405004
$ws1.Cells.Item(31, 5).Value = 91.07259052711895
$ws1.Cells.Item(31, 6).Value = 90.79266241030486
$ws1.Cells.Item(31, 7).Value = 90.58945256770775
$ws1.Cells.Item(31, 8).Value = 90.39904705689646
$ws1.Cells.Item(31, 9).Value = 90.23750468403055
$ws1.Cells.Item(31, 10).Value = 90.12691370328032
$ws1.Cells.Item(31, 11).Value = 90.00498530369269
$ws1.Cells.Item(31, 12).Value = 89.91007403799919
$ws1.Cells.Item(31, 13).Value = 89.8336811968456
$ws1.Cells.Item(31, 14).Value = 89.7721176416748
$ws1.Cells.Item(31, 15).Value = 89.70968184341267
$ws1.Cells.Item(31, 16).Value = 89.66420026847183

# --- Sheet "Budget & coverage": insert two new program rows for "Management of MAM" ---
$ws2.Rows.Item(80).Insert()
$ws2.Rows.Item(82).Insert()

# Row 80: Management of MAM / Coverage
$ws2.Cells.Item(80, 2).Value = "Management of MAM"
$ws2.Cells.Item(80, 3).Value = "Coverage"
$ws2.Cells.Item(80, 4).Value = 0
for ($col = 5; $col -le 17; $col++) {
    $ws2.Cells.Item(80, $col).Value = 0.9500000000000001
}

# Row 82: Management of MAM / Budget
$ws2.Cells.Item(82, 2).Value = "Management of MAM"
$ws2.Cells.Item(82, 3).Value = "Budget"
$ws2.Cells.Item(82, 4).Value = 0
for ($col = 5; $col -le 17; $col++) {
    $ws2.Cells.Item(82, $col).Value = 18045593.09756667
}

Write-Host "edit applied"
